# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
# A new price record (week of 2021-08-30, serial 44438) is inserted as a
# new data row right after the existing row 122, pushing the historical
# rows 123:202 down to 124:203.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 123; everything from old row 123 downward
# shifts down by one (old row 202 becomes new row 203).
$ws.Rows.Item(123).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(123, 1).Value  = 3
$ws.Cells.Item(123, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(123, 3).Value  = "Coquimbo"
$ws.Cells.Item(123, 4).Value  = 44438
$ws.Cells.Item(123, 5).Value  = 5
$ws.Cells.Item(123, 6).Value  = 100112003
$ws.Cells.Item(123, 7).Value  = "Ajo"
$ws.Cells.Item(123, 8).Value  = "Chino"
$ws.Cells.Item(123, 9).Value  = "Primera"
$ws.Cells.Item(123, 10).Value = 105
$ws.Cells.Item(123, 11).Value = 13500
$ws.Cells.Item(123, 12).Value = 14000
$ws.Cells.Item(123, 13).Value = 13738
$ws.Cells.Item(123, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(123, 15).Value = "China"
$ws.Cells.Item(123, 16).Value = 1374
$ws.Cells.Item(123, 17).Value = 10
$ws.Cells.Item(123, 18).Value = "Hortaliza"

# Keep the date column's display format consistent with the rest of the
# sheet (yyyy-mm-dd hh:mm:ss style already used by column D).
$ws.Cells.Item(123, 4).NumberFormat = $ws.Cells.Item(124, 4).NumberFormat
